$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Phased Array (Pivoting) - reduce relay mass factor from 2 to 1
$ws.Range("G4").Value = 1

# Add RELAY transmission-range data in new columns J (label) / K (range),
# entering the rows in the same order the shared-string table records them
# (row 7 up to row 3) so the underlying shared-string indices line up.
$ws.Range("J7").Value = "RELAY"
$ws.Range("K7").Value = "'1.0e+11"

$ws.Range("J6").Value = "RELAY"
$ws.Range("K6").Value = "'2.0e+11"

$ws.Range("J5").Value = "RELAY"
$ws.Range("K5").Value = "'2.0e+11"

$ws.Range("J4").Value = "RELAY"
$ws.Range("K4").Value = "'2.0e+11"

$ws.Range("J3").Value = "RELAY"
$ws.Range("K3").Value = "'1.0e+12"

# The leading apostrophes force the numeric-looking strings to be stored as
# text; strip the resulting "quote prefix" cell formatting so the cells end
# up with plain text styling (no explicit style index), matching a normal
# text entry.
$ws.Range("J3:K7").ClearFormats()

$ws.Range("M12").Select()
